$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 36) since the new data has one fewer row
$ws.Rows.Item(36).Delete()

# Update rows 2-35 with shifted values (A,B,C columns) and new D column predictions
$ws.Cells.Item(2, 1).Value = 45847.60416666666
$ws.Cells.Item(2, 2).Value = 734.655029296875
$ws.Cells.Item(2, 3).Value = 732.3300170898438
$ws.Cells.Item(2, 4).Value = 754.0988621238771
$ws.Cells.Item(3, 1).Value = 45847.64583333334
$ws.Cells.Item(3, 2).Value = 735.5700073242188
$ws.Cells.Item(3, 3).Value = 734.655029296875
$ws.Cells.Item(3, 4).Value = 727.9365735283826
$ws.Cells.Item(4, 1).Value = 45847.6875
$ws.Cells.Item(4, 2).Value = 733.77001953125
$ws.Cells.Item(4, 3).Value = 735.5700073242188
$ws.Cells.Item(4, 4).Value = 727.8367495576558
$ws.Cells.Item(5, 1).Value = 45847.72916666666
$ws.Cells.Item(5, 2).Value = 734.6599731445312
$ws.Cells.Item(5, 3).Value = 733.77001953125
$ws.Cells.Item(5, 4).Value = 733.6955798501714
$ws.Cells.Item(6, 1).Value = 45847.77083333334
$ws.Cells.Item(6, 2).Value = 735.4600219726562
$ws.Cells.Item(6, 3).Value = 734.6599731445312
$ws.Cells.Item(6, 4).Value = 751.2914037711657
$ws.Cells.Item(7, 1).Value = 45847.8125
$ws.Cells.Item(7, 2).Value = 732.7999877929688
$ws.Cells.Item(7, 3).Value = 735.4600219726562
$ws.Cells.Item(7, 4).Value = 743.311625095133
$ws.Cells.Item(8, 1).Value = 45848.5625
$ws.Cells.Item(8, 2).Value = 720.75
$ws.Cells.Item(8, 3).Value = 732.7999877929688
$ws.Cells.Item(8, 4).Value = 766.2710048749345
$ws.Cells.Item(9, 1).Value = 45848.60416666666
$ws.Cells.Item(9, 2).Value = 726.9299926757812
$ws.Cells.Item(9, 3).Value = 720.75
$ws.Cells.Item(9, 4).Value = 711.126822243416
$ws.Cells.Item(10, 1).Value = 45848.64583333334
$ws.Cells.Item(10, 2).Value = 725.964599609375
$ws.Cells.Item(10, 3).Value = 726.9299926757812
$ws.Cells.Item(10, 4).Value = 740.8774109639767
$ws.Cells.Item(11, 1).Value = 45848.6875
$ws.Cells.Item(11, 2).Value = 725.534423828125
$ws.Cells.Item(11, 3).Value = 725.964599609375
$ws.Cells.Item(11, 4).Value = 723.7343102242723
$ws.Cells.Item(12, 1).Value = 45848.72916666666
$ws.Cells.Item(12, 2).Value = 724.77001953125
$ws.Cells.Item(12, 3).Value = 725.534423828125
$ws.Cells.Item(12, 4).Value = 731.8662334862676
$ws.Cells.Item(13, 1).Value = 45848.77083333334
$ws.Cells.Item(13, 2).Value = 726.4349975585938
$ws.Cells.Item(13, 3).Value = 724.77001953125
$ws.Cells.Item(13, 4).Value = 724.101005261868
$ws.Cells.Item(14, 1).Value = 45848.8125
$ws.Cells.Item(14, 2).Value = 727.4600219726562
$ws.Cells.Item(14, 3).Value = 726.4349975585938
$ws.Cells.Item(14, 4).Value = 731.8377511133006
$ws.Cells.Item(15, 1).Value = 45849.5625
$ws.Cells.Item(15, 2).Value = 716.6199951171875
$ws.Cells.Item(15, 3).Value = 727.4600219726562
$ws.Cells.Item(15, 4).Value = 710.7900175527629
$ws.Cells.Item(16, 1).Value = 45849.60416666666
$ws.Cells.Item(16, 2).Value = 716.260009765625
$ws.Cells.Item(16, 3).Value = 716.6199951171875
$ws.Cells.Item(16, 4).Value = 729.2941470104502
$ws.Cells.Item(17, 1).Value = 45849.64583333334
$ws.Cells.Item(17, 2).Value = 720.844970703125
$ws.Cells.Item(17, 3).Value = 716.260009765625
$ws.Cells.Item(17, 4).Value = 713.4643681845046
$ws.Cells.Item(18, 1).Value = 45849.6875
$ws.Cells.Item(18, 2).Value = 721.1199951171875
$ws.Cells.Item(18, 3).Value = 720.844970703125
$ws.Cells.Item(18, 4).Value = 728.0947516097694
$ws.Cells.Item(19, 1).Value = 45849.72916666666
$ws.Cells.Item(19, 2).Value = 718.4299926757812
$ws.Cells.Item(19, 3).Value = 721.1199951171875
$ws.Cells.Item(19, 4).Value = 704.1389065990986
$ws.Cells.Item(20, 1).Value = 45849.77083333334
$ws.Cells.Item(20, 2).Value = 718.6699829101562
$ws.Cells.Item(20, 3).Value = 718.4299926757812
$ws.Cells.Item(20, 4).Value = 705.8638752329025
$ws.Cells.Item(21, 1).Value = 45849.8125
$ws.Cells.Item(21, 2).Value = 717.6099853515625
$ws.Cells.Item(21, 3).Value = 718.6699829101562
$ws.Cells.Item(21, 4).Value = 731.6788789652728
$ws.Cells.Item(22, 1).Value = 45852.5625
$ws.Cells.Item(22, 2).Value = 720.1599731445312
$ws.Cells.Item(22, 3).Value = 717.6099853515625
$ws.Cells.Item(22, 4).Value = 741.4714494469744
$ws.Cells.Item(23, 1).Value = 45852.60416666666
$ws.Cells.Item(23, 2).Value = 726.5349731445312
$ws.Cells.Item(23, 3).Value = 720.1599731445312
$ws.Cells.Item(23, 4).Value = 742.2464233412775
$ws.Cells.Item(24, 1).Value = 45852.64583333334
$ws.Cells.Item(24, 2).Value = 724.7462158203125
$ws.Cells.Item(24, 3).Value = 726.5349731445312
$ws.Cells.Item(24, 4).Value = 706.130796701729
$ws.Cells.Item(25, 1).Value = 45852.6875
$ws.Cells.Item(25, 2).Value = 724.6400146484375
$ws.Cells.Item(25, 3).Value = 724.7462158203125
$ws.Cells.Item(25, 4).Value = 727.2693403715485
$ws.Cells.Item(26, 1).Value = 45852.72916666666
$ws.Cells.Item(26, 2).Value = 724.3115234375
$ws.Cells.Item(26, 3).Value = 724.6400146484375
$ws.Cells.Item(26, 4).Value = 736.2447130184991
$ws.Cells.Item(27, 1).Value = 45852.77083333334
$ws.Cells.Item(27, 2).Value = 721.1199951171875
$ws.Cells.Item(27, 3).Value = 724.3115234375
$ws.Cells.Item(27, 4).Value = 740.2331053887306
$ws.Cells.Item(28, 1).Value = 45852.8125
$ws.Cells.Item(28, 2).Value = 720.8599853515625
$ws.Cells.Item(28, 3).Value = 721.1199951171875
$ws.Cells.Item(28, 4).Value = 697.835513657257
$ws.Cells.Item(29, 1).Value = 45853.5625
$ws.Cells.Item(29, 2).Value = 716.3099975585938
$ws.Cells.Item(29, 3).Value = 720.8599853515625
$ws.Cells.Item(29, 4).Value = 701.8181185634182
$ws.Cells.Item(30, 1).Value = 45853.60416666666
$ws.Cells.Item(30, 2).Value = 718
$ws.Cells.Item(30, 3).Value = 716.3099975585938
$ws.Cells.Item(30, 4).Value = 723.4316443202579
$ws.Cells.Item(31, 1).Value = 45853.64583333334
$ws.Cells.Item(31, 2).Value = 718.6420288085938
$ws.Cells.Item(31, 3).Value = 718
$ws.Cells.Item(31, 4).Value = 702.317226393766
$ws.Cells.Item(32, 1).Value = 45853.6875
$ws.Cells.Item(32, 2).Value = 715.2000122070312
$ws.Cells.Item(32, 3).Value = 718.6420288085938
$ws.Cells.Item(32, 4).Value = 721.6436598057328
$ws.Cells.Item(33, 1).Value = 45853.72916666666
$ws.Cells.Item(33, 2).Value = 716.5399780273438
$ws.Cells.Item(33, 3).Value = 715.2000122070312
$ws.Cells.Item(33, 4).Value = 724.1524305118799
$ws.Cells.Item(34, 1).Value = 45853.77083333334
$ws.Cells.Item(34, 2).Value = 713.1500244140625
$ws.Cells.Item(34, 3).Value = 716.5399780273438
$ws.Cells.Item(34, 4).Value = 728.1056842747284
$ws.Cells.Item(35, 1).Value = 45853.8125
$ws.Cells.Item(35, 2).Value = 710.1599731445312
$ws.Cells.Item(35, 3).Value = 713.1500244140625
$ws.Cells.Item(35, 4).Value = 747.0030100794801
